$wb = $excel.ActiveWorkbook

# Reference sheet already laid out like a quarterly "fund detail" sheet
# (bold/centered/bordered header row spanning B:H + bordered index column
# A) - used as a formatting donor for the new "2022-Q1" sheet.
$styleSrc = $wb.Worksheets.Item("2021-Q4")

# ------------------------------------------------------------------
# Step 1: duplicate the existing "总计" sheet to the end of the workbook.
# The duplicate (with correct sheetPr/pageMargins/styles already baked
# in) becomes the new "总计" sheet, while the original part is renamed
# to "2022-Q1" and repopulated with fund-detail data - this mirrors how
# the workbook was actually authored (the old "总计" sheetId/part was
# reused for "2022-Q1", and a fresh part was appended for "总计").
# ------------------------------------------------------------------
$totalOrig = $wb.Worksheets.Item("总计")
$totalOrig.Copy([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$totalSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$totalOrig.Name = "2022-Q1"
$totalSheet.Name = "总计"

$fundSheet = $totalOrig
$fundSheet.Cells.Clear()

# Header row
$fundSheet.Range("B1").Value = "基金代码"
$fundSheet.Range("C1").Value = "基金名称"
$fundSheet.Range("D1").Value = "基金规模"
$fundSheet.Range("E1").Value = "股票总仓位"
$fundSheet.Range("F1").Value = "仓位占比"
$fundSheet.Range("G1").Value = "持有市值(亿元)"
$fundSheet.Range("H1").Value = "仓位排名"

# Data rows: index(A), fund code(B), fund name(C), fund scale(D),
# total stock position(E), position ratio(F), held value(G), rank(H).
$fundData = @(
    @("006345", "景顺长城集英成长两年定期开放混合", "43.98", "88.47", "3.73", "1.6405", 9),
    @("005395", "泓德臻远回报灵活配置混合", "33.94", "93.62", "4.70", "1.5952", 8),
    @("010864", "泓德卓远混合A", "39.08", "91.88", "3.41", "1.3326", 9),
    @("010865", "泓德卓远混合C", "12.25", "91.88", "3.41", "0.4177", 9),
    @("006768", "华安沪港深优选混合", "0.84", "93.09", "3.82", "0.0321", 9)
)

$r = 2
$idx = 0
foreach ($row in $fundData) {
    $fundSheet.Range("A$r").Value = $idx
    # Fund code / scale / position% / ratio / value are kept as literal
    # text (leading zeros, fixed decimal formatting) instead of being
    # auto-converted to numbers by Excel.
    $fundSheet.Range("B$r").Value = "'" + $row[0]
    $fundSheet.Range("B$r").Style = "Normal"
    $fundSheet.Range("C$r").Value = $row[1]
    $fundSheet.Range("D$r").Value = "'" + $row[2]
    $fundSheet.Range("D$r").Style = "Normal"
    $fundSheet.Range("E$r").Value = "'" + $row[3]
    $fundSheet.Range("E$r").Style = "Normal"
    $fundSheet.Range("F$r").Value = "'" + $row[4]
    $fundSheet.Range("F$r").Style = "Normal"
    $fundSheet.Range("G$r").Value = "'" + $row[5]
    $fundSheet.Range("G$r").Style = "Normal"
    $fundSheet.Range("H$r").Value = $row[6]
    $r = $r + 1
    $idx = $idx + 1
}

# Recreate the visual style used by the other quarter sheets: bold /
# centered / bordered header row, and bordered index column.
$styleSrc.Range("B1:H1").Copy()
$fundSheet.Range("B1:H1").PasteSpecial(-4122)

$styleSrc.Range("A2").Copy()
$fundSheet.Range("A2:A6").PasteSpecial(-4122)

# ------------------------------------------------------------------
# Step 2: insert a new leading data row in "总计" for 2022-Q1, pushing
# the existing quarters down by one row.
# ------------------------------------------------------------------
$totalSheet.Rows("2:2").Insert(-4121)

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "'2022-Q1"
$totalSheet.Range("B2").Style = "Normal"
$totalSheet.Range("C2").Value = 5
$totalSheet.Range("C2").Style = "Normal"
$totalSheet.Range("D2").Value = 5.02
$totalSheet.Range("D2").Style = "Normal"

$styleSrc2 = $totalSheet
# Restyle the newly inserted row to match the rest of the table (index
# column + plain data cells with no special formatting).
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)

# Renumber the index column (A) for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4
$totalSheet.Range("A7").Value = 5

$wb.Worksheets.Item("2020-Q4").Select()
